$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Action1")

$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "test"
$ws.Range("B2").Value = "test123"
$ws.Range("B3").ClearContents()
$ws.Range("E4").Select() | Out-Null
